$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.698.48'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').Value = '1.968.67'
$ws.Range('E3').Value = '  +1.83%  '
$ws.Range('E4').Value = '  -0.01%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.23'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +0.80%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.617'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +1.28%  '
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.54'
$ws.Range('D7').Style = $origStyle
$ws.Range('E7').Value = '  +3.52%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +2.80%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0806'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -3.24%  '
$ws.Range('E11').Value = '  +0.07%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.05'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +3.44%  '
$ws.Range('D13').Value = '2.256.17'
$ws.Range('E13').Value = '  +1.82%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.824'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +1.11%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.79'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +3.12%  '
$ws.Range('E16').Value = '  +1.53%  '
$ws.Range('D17').Value = '1.966.82'
$ws.Range('E17').Value = '  +1.94%  '
$ws.Range('D18').Value = '36.573.50'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('E21').Value = '  +2.64%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.24'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E24').Value = '  -1.85%  '
$ws.Range('E25').Value = '  +2.26%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.38'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +1.48%  '
$ws.Range('E27').Value = '  -0.80%  '
$ws.Range('E28').Value = '  +11.78%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.36'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +0.33%  '
$ws.Range('E30').Value = '  +1.74%  '
$ws.Range('E31').Value = '  -1.19%  '
$ws.Range('E32').Value = '  +1.06%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0619'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  -1.02%  '
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('B36').Value = 'THORChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.08'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +1.63%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.40'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  +17.34%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.24'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +4.88%  '
$ws.Range('E39').Value = '  -0.74%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0996'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +2.95%  '
$ws.Range('E41').Value = '  +1.74%  '
$ws.Range('E42').Value = '  +2.54%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.16'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  -0.16%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.08'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  +2.78%  '
$ws.Range('D45').Value = '1.368.36'
$ws.Range('E45').Value = '  +2.10%  '
$ws.Range('E46').Value = '  +1.35%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '87.53'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  +0.47%  '
$ws.Range('E48').Value = '  -0.30%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.84'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +1.15%  '
$ws.Range('D50').Value = '2.147.14'
$ws.Range('E50').Value = '  +1.85%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.34'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  -4.27%  '
